$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3947.25
$ws.Range("J51").Value = 3947.25
$ws.Range("L51").Value = 3947.25
$ws.Range("N51").Value = -4915.25
$ws.Range("H58").Value = 2369.8
$ws.Range("I58").Value = 539.6
$ws.Range("J58").Value = 4200
$ws.Range("K58").Value = 1618.8
$ws.Range("L58").Value = 12600
$ws.Range("M58").Value = -1468.8
$ws.Range("N58").Value = -12900
$ws.Range("H70").Value = 3420.5789
$ws.Range("I70").Value = 3827.0908
$ws.Range("J70").Value = 2861.625
$ws.Range("K70").Value = 11481.2724
$ws.Range("L70").Value = 8584.875
$ws.Range("M70").Value = -11211.2724
$ws.Range("N70").Value = -9124.875
$ws.Range("H73").Value = 3420.5789
$ws.Range("I73").Value = 3827.0908
$ws.Range("J73").Value = 2861.625
$ws.Range("K73").Value = 11481.2724
$ws.Range("L73").Value = 8584.875
$ws.Range("M73").Value = -10545.2724
$ws.Range("N73").Value = -10456.875
$ws.Range("H76").Value = 4950
$ws.Range("I76").Value = 4900
$ws.Range("K76").Value = 4900
$ws.Range("M76").Value = -4585
$ws.Range("H79").Value = 4950
$ws.Range("I79").Value = 4900
$ws.Range("K79").Value = 4900
$ws.Range("M79").Value = -3808
$ws.Range("H95").Value = 29999
$ws.Range("J95").Value = 29999
$ws.Range("L95").Value = 29999
$ws.Range("N95").Value = -35491
$ws.Range("H100").Value = 4815415.5
$ws.Range("I100").Value = 6477.3076
$ws.Range("J100").Value = 15234781
$ws.Range("K100").Value = 6477.3076
$ws.Range("L100").Value = 15234781
$ws.Range("M100").Value = -5936.3076
$ws.Range("N100").Value = -15235863
$ws.Range("H137").Value = 560293.75
$ws.Range("I137").Value = 658649.8
$ws.Range("J137").Value = 35728
$ws.Range("K137").Value = 1975949.4
$ws.Range("L137").Value = 107184
$ws.Range("M137").Value = -1973399.4
$ws.Range("N137").Value = -112284
$ws.Range("H138").Value = 4645.656
$ws.Range("I138").Value = 1321.4166
$ws.Range("J138").Value = 5138.1357
$ws.Range("K138").Value = 3964.2498
$ws.Range("L138").Value = 15414.4071
$ws.Range("M138").Value = 1175.7502
$ws.Range("N138").Value = -25694.4071

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1281.5714
$ws.Range("I5").Value = 328.5
$ws.Range("J5").Value = 7000
$ws.Range("K5").Value = 328.5
$ws.Range("L5").Value = 7000
$ws.Range("M5").Value = -216.5
$ws.Range("N5").Value = -7224
$ws.Range("H16").Value = 5835.6665
$ws.Range("I16").Value = 7250
$ws.Range("J16").Value = 3007
$ws.Range("K16").Value = 7250
$ws.Range("L16").Value = 3007
$ws.Range("M16").Value = -6963
$ws.Range("N16").Value = -3581
$ws.Range("H32").Value = 2932.8462
$ws.Range("I32").Value = 3008.5
$ws.Range("K32").Value = 3008.5
$ws.Range("M32").Value = -2721.5
$ws.Range("H61").Value = 11272.714
$ws.Range("I61").Value = 19299
$ws.Range("K61").Value = 19299
$ws.Range("M61").Value = -19087
$ws.Range("H88").Value = 1444.3334
$ws.Range("J88").Value = 1921.7142
$ws.Range("L88").Value = 1921.7142
$ws.Range("N88").Value = -2733.7142
$ws.Range("H91").Value = 1444.3334
$ws.Range("J91").Value = 1921.7142
$ws.Range("L91").Value = 1921.7142
$ws.Range("N91").Value = -4729.7142
$ws.Range("H97").Value = 5480.231
$ws.Range("I97").Value = 8140.1875
$ws.Range("J97").Value = 1224.3
$ws.Range("K97").Value = 8140.1875
$ws.Range("L97").Value = 1224.3
$ws.Range("M97").Value = -7644.1875
$ws.Range("N97").Value = -2216.3
$ws.Range("H110").Value = 2000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 2000
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -6090
$ws.Range("H122").Value = 672179.6
$ws.Range("I122").Value = 4052.0667
$ws.Range("K122").Value = 12156.2001
$ws.Range("M122").Value = -9706.2001
$ws.Range("H132").Value = 3136.96
$ws.Range("I132").Value = 2193.3684
$ws.Range("K132").Value = 6580.1052
$ws.Range("M132").Value = -4050.1052
$ws.Range("H136").Value = 11272.714
$ws.Range("I136").Value = 19299
$ws.Range("K136").Value = 57897
$ws.Range("M136").Value = -55347

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1281.5714
$ws.Range("I4").Value = 328.5
$ws.Range("J4").Value = 7000
$ws.Range("K4").Value = 328.5
$ws.Range("L4").Value = 7000
$ws.Range("M4").Value = -213.5
$ws.Range("N4").Value = -7230
$ws.Range("H94").Value = 2232.04
$ws.Range("J94").Value = 3884.8
$ws.Range("L94").Value = 3884.8
$ws.Range("N94").Value = -4786.8
$ws.Range("H105").Value = 2772.5
$ws.Range("I105").Value = 2295.3333
$ws.Range("J105").Value = 3488.25
$ws.Range("K105").Value = 2295.3333
$ws.Range("L105").Value = 3488.25
$ws.Range("M105").Value = -548.3332999999998
$ws.Range("N105").Value = -6982.25
$ws.Range("H107").Value = 5256.6
$ws.Range("I107").Value = 5714.636
$ws.Range("J107").Value = 3997
$ws.Range("K107").Value = 5714.636
$ws.Range("L107").Value = 3997
$ws.Range("M107").Value = -3794.636
$ws.Range("N107").Value = -7837
$ws.Range("H134").Value = 5467.3076
$ws.Range("I134").Value = 5287.6
$ws.Range("K134").Value = 15862.8
$ws.Range("M134").Value = -13327.8

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 8647.27
$ws.Range("I86").Value = 8339.764999999999
$ws.Range("K86").Value = 8339.764999999999
$ws.Range("M86").Value = -7216.764999999999
$ws.Range("H89").Value = 8647.27
$ws.Range("I89").Value = 8339.764999999999
$ws.Range("K89").Value = 41698.825
$ws.Range("M89").Value = -36082.825

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 262.95
$ws.Range("J2").Value = 282
$ws.Range("L2").Value = 1692
$ws.Range("N2").Value = -1918
$ws.Range("H5").Value = 358701.97
$ws.Range("I5").Value = 1821.2142
$ws.Range("K5").Value = 5463.642599999999
$ws.Range("M5").Value = -5351.642599999999
$ws.Range("H23").Value = 253
$ws.Range("I23").Value = 56
$ws.Range("J23").Value = 450
$ws.Range("K23").Value = 168
$ws.Range("L23").Value = 1350
$ws.Range("M23").Value = 67
$ws.Range("N23").Value = -1820
$ws.Range("H34").Value = 1925264.9
$ws.Range("I34").Value = 2780049.2
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 8340147.600000001
$ws.Range("L34").Value = 6000
$ws.Range("M34").Value = -8340063.600000001
$ws.Range("N34").Value = -6168
$ws.Range("H48").Value = 3274.25
$ws.Range("J48").Value = 9600
$ws.Range("L48").Value = 28800
$ws.Range("N48").Value = -29300
$ws.Range("H55").Value = 10040.366
$ws.Range("I55").Value = 1225.75
$ws.Range("J55").Value = 11396.462
$ws.Range("K55").Value = 3677.25
$ws.Range("L55").Value = 34189.386
$ws.Range("M55").Value = -3500.25
$ws.Range("N55").Value = -34543.386
$ws.Range("H80").Value = 100460.92
$ws.Range("I80").Value = 3499.5
$ws.Range("K80").Value = 10498.5
$ws.Range("M80").Value = -9562.5
$ws.Range("H83").Value = 100460.92
$ws.Range("I83").Value = 3499.5
$ws.Range("K83").Value = 31495.5
$ws.Range("M83").Value = -26815.5
$ws.Range("H108").Value = 2338.2727
$ws.Range("I108").Value = 1472.1
$ws.Range("K108").Value = 4416.299999999999
$ws.Range("M108").Value = -1536.299999999999
$ws.Range("H122").Value = 2024.2188
$ws.Range("I122").Value = 956
$ws.Range("J122").Value = 2380.2917
$ws.Range("K122").Value = 8604
$ws.Range("L122").Value = 21422.6253
$ws.Range("M122").Value = -6154
$ws.Range("N122").Value = -26322.6253
$ws.Range("H135").Value = 358701.97
$ws.Range("I135").Value = 1821.2142
$ws.Range("K135").Value = 16390.9278
$ws.Range("M135").Value = -13855.9278

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 16464.05
$ws.Range("I126").Value = 33999.8
$ws.Range("J126").Value = 10618.8
$ws.Range("K126").Value = 101999.4
$ws.Range("L126").Value = 31856.4
$ws.Range("M126").Value = -99529.40000000001
$ws.Range("N126").Value = -36796.39999999999
$ws.Range("H132").Value = 3356.5925
$ws.Range("I132").Value = 2407.2
$ws.Range("J132").Value = 6069.143
$ws.Range("K132").Value = 7221.599999999999
$ws.Range("L132").Value = 18207.429
$ws.Range("M132").Value = -4691.599999999999
$ws.Range("N132").Value = -23267.429

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2317.4167
$ws.Range("I46").Value = 1769.4166
$ws.Range("K46").Value = 1769.4166
$ws.Range("M46").Value = -1581.4166
$ws.Range("H132").Value = 1449942.4
$ws.Range("I132").Value = 1619194.6
$ws.Range("J132").Value = 11298.5
$ws.Range("K132").Value = 4857583.800000001
$ws.Range("L132").Value = 33895.5
$ws.Range("M132").Value = -4855053.800000001
$ws.Range("N132").Value = -38955.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 79986
$ws.Range("J46").Value = 79986
$ws.Range("L46").Value = 79986
$ws.Range("N46").Value = -80448
$ws.Range("H92").Value = 275005000
$ws.Range("I92").Value = 10000
$ws.Range("K92").Value = 10000
$ws.Range("M92").Value = -7504
$ws.Range("H126").Value = 29596.688
$ws.Range("J126").Value = 7874.4
$ws.Range("L126").Value = 23623.2
$ws.Range("N126").Value = -28563.2
$ws.Range("H132").Value = 11330.128
$ws.Range("I132").Value = 13535.032
$ws.Range("K132").Value = 40605.096
$ws.Range("M132").Value = -38075.096
$ws.Range("H134").Value = 79986
$ws.Range("J134").Value = 79986
$ws.Range("L134").Value = 239958
$ws.Range("N134").Value = -245028

